# NIT-9007567820.xlsx — "Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta"
#
# The worker rows (16-19) on Hoja1 get re-sorted (grouped by period instead of
# by worker) and Luis Alberto Hernandez Julio's "Salario Basico" is updated
# from 737717 to 781242. Jaime Luis Leones Villa's salary is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: LUIS ALBERTO HERNANDEZ JULIO - periodo 1707 - salario actualizado
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1143340528"
$ws.Range("D16").Value = "LUIS ALBERTO HERNANDEZ JULIO"
$ws.Range("E16").Value = "1707"
$ws.Range("F16").Value = 14754
$ws.Range("G16").Value = 781242

# Row 17: JAIME LUIS LEONES VILLA - periodo 1707
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "73230640"
$ws.Range("D17").Value = "JAIME LUIS LEONES VILLA"
$ws.Range("E17").Value = "1707"
$ws.Range("F17").Value = 14754
$ws.Range("G17").Value = 737717

# Row 18: LUIS ALBERTO HERNANDEZ JULIO - periodo 1708 - salario actualizado
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1143340528"
$ws.Range("D18").Value = "LUIS ALBERTO HERNANDEZ JULIO"
$ws.Range("E18").Value = "1708"
$ws.Range("F18").Value = 14754
$ws.Range("G18").Value = 781242

# Row 19: JAIME LUIS LEONES VILLA - periodo 1708
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "73230640"
$ws.Range("D19").Value = "JAIME LUIS LEONES VILLA"
$ws.Range("E19").Value = "1708"
$ws.Range("F19").Value = 14754
$ws.Range("G19").Value = 737717
